$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-16 (A2:A16)
$valuesTop = @(
    3.977959584536734,
    7.622334346345951,
    8.073032840776023,
    9.379189874814472,
    2.736651432341773,
    9.799971428978751,
    6.595681556606564,
    5.560753044758911,
    7.876287900559333,
    8.181062382652982,
    1.446861457488097,
    5.792735607280207,
    4.063784253847302,
    2.792715961972817,
    2.431485653114777
)

for ($i = 0; $i -lt $valuesTop.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $valuesTop[$i]
}

# Rows 17-86 all take the same constant value
$constValue = 5.813051703908144
for ($row = 17; $row -le 86; $row++) {
    $ws.Cells.Item($row, 1).Value = $constValue
}
